# Update the "想去人数" (number of people interested) figures that changed
# between scrapes, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F5").Value = 11561
    $ws.Range("F14").Value = 51
    $ws.Range("F18").Value = 1324
    $ws.Range("F19").Value = 74
    $ws.Range("F20").Value = 900
}
